$wb = $excel.ActiveWorkbook

# Data for 2025-07-28: updated column L ("2025" year-to-date totals)
# across Citywide Totals, By Neighborhood summary, and individual neighborhood sheets.
$sheetUpdates = @(
    @("Citywide Totals", @(@("L2", 3849), @("L3", 4029), @("J4", 1870), @("L4", 994), @("L5", 235), @("L6", 3480), @("J7", 29346), @("L7", 12587))),
    @("Logan Square", @(@("L2", 43), @("L3", 37), @("L7", 141))),
    @("Austin", @(@("L2", 239), @("L6", 224), @("L7", 817))),
    @("Garfield Park", @(@("L2", 164), @("L3", 195), @("L4", 33), @("L6", 188), @("L7", 591))),
    @("Grand Crossing", @(@("L2", 138), @("L3", 152), @("L6", 126), @("L7", 460))),
    @("Woodlawn", @(@("L6", 43), @("L7", 207))),
    @("By Neighborhood", @(@("L2", 107), @("L6", 103), @("L7", 413), @("L8", 817), @("L15", 91), @("L17", 24), @("L19", 360), @("L20", 313), @("L22", 36), @("L29", 693), @("L33", 591), @("L34", 78), @("L37", 460), @("L42", 403), @("L48", 169), @("L49", 67), @("L52", 253), @("L53", 141), @("L54", 259), @("L55", 117), @("L57", 42), @("J63", 222), @("L63", 41), @("L64", 83), @("L73", 105), @("L76", 190), @("L77", 85), @("L78", 158), @("L79", 326), @("L80", 40), @("L85", 661), @("L86", 93), @("L88", 140), @("L90", 121), @("L94", 159), @("L96", 132), @("L97", 105), @("L98", 72), @("L99", 207), @("J101", 29346), @("L101", 12587))),
    @("Lincoln Park", @(@("L2", 24), @("L7", 67))),
    @("Loop", @(@("L6", 120), @("L7", 259))),
    @("Englewood", @(@("L2", 215), @("L5", 11), @("L7", 693))),
    @("Lake View", @(@("L3", 39), @("L7", 169))),
    @("Chatham", @(@("L3", 111), @("L7", 360))),
    @("River North", @(@("L2", 38), @("L3", 35), @("L7", 190))),
    @("Ashburn", @(@("L3", 29), @("L7", 103))),
    @("Humboldt Park", @(@("L2", 121), @("L3", 130), @("L7", 403))),
    @("Rogers Park", @(@("L2", 45), @("L7", 158))),
    @("Lower West Side", @(@("L2", 38), @("L7", 117))),
    @("West Ridge", @(@("L3", 35), @("L4", 17), @("L6", 35), @("L7", 132))),
    @("Roseland", @(@("L6", 66), @("L7", 326))),
    @("Near South Side", @(@("L2", 27), @("L7", 83))),
    @("Chicago Lawn", @(@("L2", 98), @("L3", 103), @("L7", 313))),
    @("Burnside", @(@("L3", 8), @("L7", 24))),
    @("Auburn Gresham", @(@("L2", 136), @("L3", 127), @("L7", 413))),
    @("Garfield Ridge", @(@("L2", 23), @("L3", 20), @("L7", 78))),
    @("West Loop", @(@("L6", 62), @("L7", 159))),
    @("Brighton Park", @(@("L3", 28), @("L7", 91))),
    @("Wicker Park", @(@("L3", 11), @("L7", 72))),
    @("Portage Park", @(@("L3", 31), @("L7", 105))),
    @("Albany Park", @(@("L2", 36), @("L7", 107))),
    @("West Town", @(@("L6", 54), @("L7", 105))),
    @("United Center", @(@("L2", 37), @("L7", 140))),
    @("Streeterville", @(@("L6", 11), @("L7", 93))),
    @("Washington Heights", @(@("L2", 41), @("L7", 121))),
    @("Mckinley Park", @(@("L2", 15), @("L7", 42))),
    @("South Shore", @(@("L2", 196), @("L3", 268), @("L7", 661))),
    @("Clearing", @(@("L4", 3), @("L7", 36))),
    @("Riverdale", @(@("L2", 29), @("L7", 85))),
    @("Rush & Division", @(@("L3", 9), @("L7", 40))),
    @("Little Village", @(@("L2", 87), @("L7", 253)))
)

foreach ($sheetUpdate in $sheetUpdates) {
    $sheetName = $sheetUpdate[0]
    $cellPairs = $sheetUpdate[1]
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $cellPairs) {
        $ws.Range($pair[0]).Value = $pair[1]
    }
}
